$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save" - copy formatting from the adjacent "sum" header (G1)
# so it gets the same bold/centered/bordered style used by the other header cells.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# New data column values for rows 2 and 3 (plain numeric, no special style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0

$excel.CutCopyMode = 0
